# Update product name in A3 from "Smartwatch " to "Smart watch "
# and move the selection/active cell to A3 (cosmetic, matches the saved view state).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Smart watch "

$ws.Range("A3").Select()
